$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Status text constants (reuse existing shared strings)
$noRequisites = "Исполнитель не оставил свои реквизиты"
$closed = "Закрыт"
$notClosed = "Не закрыт"

# Row 2
$ws.Cells.Item(2, 1).Value = 23
$ws.Cells.Item(2, 2).Value = 185404885
$ws.Cells.Item(2, 3).Value = 236322856
$ws.Cells.Item(2, 4).Value = 45041.86893779451
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 45041.50833333333
$ws.Cells.Item(2, 7).Value = 1212
$ws.Cells.Item(2, 8).Value = 10
$ws.Cells.Item(2, 9).Value = 121
$ws.Cells.Item(2, 10).Value = 1090
$ws.Cells.Item(2, 11).Value = $noRequisites
$ws.Cells.Item(2, 12).Value = $closed

# Row 3
$ws.Cells.Item(3, 1).Value = 24
$ws.Cells.Item(3, 2).Value = 185404885
$ws.Cells.Item(3, 3).Value = 236322856
$ws.Cells.Item(3, 4).Value = 45041.883623359245
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 45041.50833333333
$ws.Cells.Item(3, 7).Value = 1200
$ws.Cells.Item(3, 8).Value = 10
$ws.Cells.Item(3, 9).Value = 120
$ws.Cells.Item(3, 10).Value = 1080
$ws.Cells.Item(3, 11).Value = $noRequisites
$ws.Cells.Item(3, 12).Value = $closed

# Row 4
$ws.Cells.Item(4, 1).Value = 25
$ws.Cells.Item(4, 2).Value = 185404885
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 45041.886526336224
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 45041.50833333333
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = $noRequisites
$ws.Cells.Item(4, 12).Value = $notClosed

# Row 5
$ws.Cells.Item(5, 1).Value = 26
$ws.Cells.Item(5, 2).Value = 185404885
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = 45042.84731792426
$ws.Cells.Item(5, 5).Value = 4
$ws.Cells.Item(5, 6).Value = 45049.50833333333
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = $noRequisites
$ws.Cells.Item(5, 12).Value = $notClosed

# Row 6
$ws.Cells.Item(6, 1).Value = 27
$ws.Cells.Item(6, 2).Value = 185404885
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 45042.86098662064
$ws.Cells.Item(6, 5).Value = 5
$ws.Cells.Item(6, 6).Value = 45043.50833333333
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = $noRequisites
$ws.Cells.Item(6, 12).Value = $notClosed

# Row 7
$ws.Cells.Item(7, 1).Value = 28
$ws.Cells.Item(7, 2).Value = 185404885
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 45042.86229964735
$ws.Cells.Item(7, 5).Value = 12
$ws.Cells.Item(7, 6).Value = 45056.50833333333
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = $noRequisites
$ws.Cells.Item(7, 12).Value = $notClosed

# Row 8
$ws.Cells.Item(8, 1).Value = 22
$ws.Cells.Item(8, 2).Value = 185404885
$ws.Cells.Item(8, 3).Value = 236322856
$ws.Cells.Item(8, 4).Value = 45041.86554414516
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 45041.50833333333
$ws.Cells.Item(8, 7).Value = 234
$ws.Cells.Item(8, 8).Value = 10
$ws.Cells.Item(8, 9).Value = 23
$ws.Cells.Item(8, 10).Value = 210
$ws.Cells.Item(8, 11).Value = $noRequisites
$ws.Cells.Item(8, 12).Value = $closed


# Fix number formatting (date style) for newly added rows 7 and 8 (D and F columns)
$ws.Range("D2").Copy()
$ws.Range("D7:D8").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("F7:F8").PasteSpecial(-4122)

# Fix status fill styles on column L based on target status
# Style 1 = green fill ("Закрыт"); Style 3 = red fill ("Не закрыт")
# L2 is already style 1 (green) and stays that way -> use as the green source.
# L6 is already style 3 (red) and stays that way -> use as the red source.
# (Multi-area ranges, e.g. "L3,L8", drop all but the first area when used
# with PasteSpecial here, so each target cell is pasted individually.)
$ws.Range("L2").Copy()
$ws.Range("L3").PasteSpecial(-4122)
$ws.Range("L2").Copy()
$ws.Range("L8").PasteSpecial(-4122)
$ws.Range("L6").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L6").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L6").Copy()
$ws.Range("L7").PasteSpecial(-4122)

$excel.CutCopyMode = $false

